# Atualizando o arquivo XLSX
# Update odds values (Jogos_da_Semana_FlashScore_2024-11-18.xlsx)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2
$ws.Range("H2").Value = 3.3
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 2.75
$ws.Range("L2").Value = 4.5
$ws.Range("M2").Value = 1.07
$ws.Range("N2").Value = 9
$ws.Range("O2").Value = 1.36
$ws.Range("Q2").Value = 2.2
$ws.Range("R2").Value = 1.65
$ws.Range("S2").Value = 1.5
$ws.Range("T2").Value = 2.5
$ws.Range("U2").Value = 1.95
$ws.Range("V2").Value = 1.8
$ws.Range("X2").Value = 8.5
$ws.Range("Z2").Value = 17
$ws.Range("AB2").Value = 34
$ws.Range("AC2").Value = 8
$ws.Range("AI2").Value = 19
$ws.Range("AJ2").Value = 13
$ws.Range("AL2").Value = 34
$ws.Range("AO2").Value = 11
$ws.Range("AR2").Value = 67
$ws.Range("AT2").Value = 2.5
$ws.Range("BA2").Value = 101

# Row 3
$ws.Range("G3").Value = 1.9
$ws.Range("K3").Value = 1.83
$ws.Range("M3").Value = 1.14
$ws.Range("O3").Value = 1.67
$ws.Range("R3").Value = 1.36

# Row 4
$ws.Range("G4").Value = 2.2
$ws.Range("K4").Value = 1.91
$ws.Range("M4").Value = 1.11
$ws.Range("O4").Value = 1.53
$ws.Range("P4").Value = 2.38
$ws.Range("R4").Value = 1.44

# Row 5
$ws.Range("G5").Value = 2.3
$ws.Range("K5").Value = 1.91
$ws.Range("M5").Value = 1.1
$ws.Range("O5").Value = 1.5
$ws.Range("R5").Value = 1.48

# Row 6
$ws.Range("G6").Value = 2.45
$ws.Range("H6").Value = 2.8
$ws.Range("K6").Value = 1.8
$ws.Range("M6").Value = 1.13
$ws.Range("O6").Value = 1.62
$ws.Range("R6").Value = 1.36

# Row 8
$ws.Range("G8").Value = 3.3
$ws.Range("I8").Value = 2.38
$ws.Range("L8").Value = 3.2
$ws.Range("AD8").Value = 5.5
$ws.Range("AH8").Value = 7
$ws.Range("AP8").Value = 29
$ws.Range("AZ8").Value = 51
